$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '25.818.10'
$c.Style = "Normal"
$ws.Range('E2').Value = '  -2.50%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '1.745.42'
$c.Style = "Normal"
$ws.Range('E3').Value = '  -4.98%  '
$c = $ws.Range('D4')
$c.NumberFormat = "@"
$c.Value = '0.9999'
$c.Style = "Normal"
$ws.Range('E4').Value = '  -0.10%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '237.41'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -9.05%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.9993'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -0.18%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.5052'
$c.Style = "Normal"
$ws.Range('E7').Value = '  -5.78%  '
$ws.Range('E8').Value = '  -6.63%  '
$ws.Range('E9').Value = '  -12.16%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.06159'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -10.33%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '1.743.20'
$c.Style = "Normal"
$ws.Range('E11').Value = '  -5.35%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '0.06919'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -4.25%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '15.42'
$c.Style = "Normal"
$ws.Range('E13').Value = '  -12.16%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '0.5994'
$c.Style = "Normal"
$ws.Range('E14').Value = '  -18.88%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '4.497'
$c.Style = "Normal"
$ws.Range('E15').Value = '  -9.58%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '76.91'
$c.Style = "Normal"
$ws.Range('E16').Value = '  -14.06%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '0.9994'
$c.Style = "Normal"
$ws.Range('E17').Value = '  -0.16%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '0.9995'
$c.Style = "Normal"
$ws.Range('E18').Value = '  -0.15%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '25.828.12'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -2.56%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '0.000006850'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -13.00%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '11.60'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -15.95%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '1.964.32'
$c.Style = "Normal"
$ws.Range('E22').Value = '  -5.55%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '4.039'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -11.87%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '5.201'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -12.78%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '8.132'
$c.Style = "Normal"
$ws.Range('E25').Value = '  -11.88%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '137.92'
$c.Style = "Normal"
$ws.Range('E26').Value = '  -3.43%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '1.515'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -10.02%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '14.99'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -11.55%  '
$ws.Range('E29').Value = '  -17.24%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '103.46'
$c.Style = "Normal"
$ws.Range('E30').Value = '  -6.51%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '3.766'
$c.Style = "Normal"
$ws.Range('E31').Value = '  -10.70%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '0.08105'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -8.04%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '3.469'
$c.Style = "Normal"
$ws.Range('E33').Value = '  -13.79%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '0.04516'
$c.Style = "Normal"
$ws.Range('E34').Value = '  -6.15%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '0.9988'
$c.Style = "Normal"
$ws.Range('E35').Value = '  -0.13%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '2.651'
$c.Style = "Normal"
$ws.Range('E36').Value = '  -9.38%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.9836'
$c.Style = "Normal"
$ws.Range('E37').Value = '  -12.99%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.6072'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -16.81%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '2.672'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -13.73%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.01548'
$c.Style = "Normal"
$ws.Range('E40').Value = '  -9.25%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '1.910'
$c.Style = "Normal"
$ws.Range('E41').Value = '  -16.03%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '0.9990'
$c.Style = "Normal"
$ws.Range('E42').Value = '  -0.18%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '103.08'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -4.27%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '0.3809'
$c.Style = "Normal"
$ws.Range('E44').Value = '  -19.15%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '5.088'
$c.Style = "Normal"
$ws.Range('E45').Value = '  -13.53%  '
$ws.Range('E46').Value = '  -19.14%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '0.05349'
$c.Style = "Normal"
$ws.Range('E47').Value = '  -7.45%  '
$ws.Range('E48').Value = '  -10.06%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '30.20'
$c.Style = "Normal"
$ws.Range('E49').Value = '  -12.96%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '5.915'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -19.70%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '52.58'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -12.39%  '
